# Weekly refresh of "Fruta, Feria Lagunitas de Puerto Montt - Granada" price records.
# The underlying daily price feed was re-pulled for the week, which shuffles the
# Fecha / Calidad / Volumen / Precio* / Unidad / Precio $/Kg / Kg-unidad columns
# across the existing rows (Mercado, Producto, Categoria, Origen, etc. stay fixed).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44351

# Row 4
$ws.Range("D4").Value = 44351

# Row 5
$ws.Range("D5").Value = 44313
$ws.Range("L5").Value = 'Especial'
$ws.Range("M5").Value = 100
$ws.Range("Q5").Value = '$/caja 14 kilos empedrada'
$ws.Range("S5").Value = 1250
$ws.Range("T5").Value = 14

# Row 6
$ws.Range("D6").Value = 44313
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 16000
$ws.Range("O6").Value = 16000
$ws.Range("P6").Value = 16000
$ws.Range("Q6").Value = '$/caja 14 kilos empedrada'
$ws.Range("S6").Value = 1143
$ws.Range("T6").Value = 14

# Row 7
$ws.Range("D7").Value = 44313
$ws.Range("L7").Value = 'Segunda'
$ws.Range("M7").Value = 80
$ws.Range("N7").Value = 14000
$ws.Range("O7").Value = 14000
$ws.Range("P7").Value = 14000
$ws.Range("Q7").Value = '$/caja 14 kilos empedrada'
$ws.Range("S7").Value = 1000
$ws.Range("T7").Value = 14

# Row 8
$ws.Range("D8").Value = 44334
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 200
$ws.Range("N8").Value = 14000
$ws.Range("O8").Value = 17000
$ws.Range("P8").Value = 15500
$ws.Range("S8").Value = 1033

# Row 9
$ws.Range("D9").Value = 44334
$ws.Range("L9").Value = 'Segunda'
$ws.Range("N9").Value = 14500
$ws.Range("O9").Value = 14500
$ws.Range("P9").Value = 14500
$ws.Range("Q9").Value = '$/caja 15 kilos empedrada'
$ws.Range("S9").Value = 967
$ws.Range("T9").Value = 15

# Row 10
$ws.Range("D10").Value = 44316
$ws.Range("N10").Value = 17500
$ws.Range("O10").Value = 17500
$ws.Range("P10").Value = 17500
$ws.Range("Q10").Value = '$/caja 15 kilos empedrada'
$ws.Range("S10").Value = 1167
$ws.Range("T10").Value = 15

# Row 11
$ws.Range("D11").Value = 44316
$ws.Range("M11").Value = 200
$ws.Range("O11").Value = 14500
$ws.Range("P11").Value = 14250
$ws.Range("Q11").Value = '$/caja 15 kilos empedrada'
$ws.Range("S11").Value = 950
$ws.Range("T11").Value = 15

# Row 12
$ws.Range("D12").Value = 44344
$ws.Range("N12").Value = 16000
$ws.Range("O12").Value = 16000
$ws.Range("P12").Value = 16000
$ws.Range("S12").Value = 1067

# Row 13
$ws.Range("D13").Value = 44344
$ws.Range("M13").Value = 120
$ws.Range("N13").Value = 13000
$ws.Range("O13").Value = 13500
$ws.Range("P13").Value = 13250
$ws.Range("S13").Value = 883

# Row 14
$ws.Range("D14").Value = 44298
$ws.Range("L14").Value = 'Segunda'
$ws.Range("M14").Value = 80
$ws.Range("N14").Value = 14000
$ws.Range("O14").Value = 15000
$ws.Range("P14").Value = 14500
$ws.Range("S14").Value = 967

# Row 15
$ws.Range("D15").Value = 44293
$ws.Range("L15").Value = 'Primera'
$ws.Range("M15").Value = 60
$ws.Range("O15").Value = 15000
$ws.Range("S15").Value = 967

# Row 16
$ws.Range("D16").Value = 44302

# Row 17
$ws.Range("D17").Value = 44302
$ws.Range("O17").Value = 15000
$ws.Range("S17").Value = 967

# Row 18
$ws.Range("D18").Value = 44292
$ws.Range("M18").Value = 160

# Row 19
$ws.Range("D19").Value = 44306
$ws.Range("M19").Value = 100

# Row 20
$ws.Range("D20").Value = 44306
$ws.Range("M20").Value = 200
$ws.Range("O20").Value = 14500
$ws.Range("P20").Value = 14250
$ws.Range("S20").Value = 950

# Row 21
$ws.Range("D21").Value = 44323
$ws.Range("M21").Value = 100
$ws.Range("N21").Value = 17000
$ws.Range("P21").Value = 17000
$ws.Range("S21").Value = 1133

# Row 22
$ws.Range("D22").Value = 44323
$ws.Range("N22").Value = 14000
$ws.Range("O22").Value = 14000
$ws.Range("P22").Value = 14000
$ws.Range("S22").Value = 933

# Row 23
$ws.Range("D23").Value = 44295
$ws.Range("M23").Value = 160
$ws.Range("N23").Value = 14000
$ws.Range("O23").Value = 15000
$ws.Range("P23").Value = 14500
$ws.Range("S23").Value = 967

# Row 24
$ws.Range("D24").Value = 44299
$ws.Range("L24").Value = 'Primera'
$ws.Range("M24").Value = 60
$ws.Range("N24").Value = 17500
$ws.Range("O24").Value = 17500
$ws.Range("P24").Value = 17500
$ws.Range("S24").Value = 1167

# Row 25
$ws.Range("D25").Value = 44299
$ws.Range("L25").Value = 'Segunda'
$ws.Range("M25").Value = 120
$ws.Range("N25").Value = 14000
$ws.Range("P25").Value = 14500
$ws.Range("S25").Value = 967

# Row 26
$ws.Range("D26").Value = 44336
$ws.Range("L26").Value = 'Primera'
$ws.Range("M26").Value = 60
$ws.Range("N26").Value = 17000
$ws.Range("O26").Value = 17000
$ws.Range("P26").Value = 17000
$ws.Range("S26").Value = 1133

# Row 27
$ws.Range("D27").Value = 44336
$ws.Range("L27").Value = 'Segunda'
$ws.Range("M27").Value = 120
$ws.Range("O27").Value = 14500
$ws.Range("P27").Value = 14250
$ws.Range("S27").Value = 950

# Row 28
$ws.Range("D28").Value = 44301
$ws.Range("L28").Value = 'Primera'
$ws.Range("M28").Value = 60
$ws.Range("N28").Value = 17500
$ws.Range("O28").Value = 17500
$ws.Range("P28").Value = 17500
$ws.Range("S28").Value = 1167

# Row 29
$ws.Range("D29").Value = 44301
$ws.Range("L29").Value = 'Segunda'
$ws.Range("M29").Value = 80
$ws.Range("N29").Value = 14000
$ws.Range("O29").Value = 15000
$ws.Range("P29").Value = 14500
$ws.Range("S29").Value = 967

# Row 30
$ws.Range("D30").Value = 44305
$ws.Range("L30").Value = 'Primera'
$ws.Range("M30").Value = 60
$ws.Range("N30").Value = 17500
$ws.Range("O30").Value = 17500
$ws.Range("P30").Value = 17500
$ws.Range("S30").Value = 1167

# Row 31
$ws.Range("D31").Value = 44305
$ws.Range("L31").Value = 'Segunda'
$ws.Range("M31").Value = 120
$ws.Range("N31").Value = 14000
$ws.Range("O31").Value = 15000
$ws.Range("P31").Value = 14500
$ws.Range("S31").Value = 967

# Row 32
$ws.Range("D32").Value = 44348
$ws.Range("L32").Value = 'Primera'
$ws.Range("M32").Value = 100
$ws.Range("N32").Value = 15000
$ws.Range("O32").Value = 15000
$ws.Range("P32").Value = 15000
$ws.Range("S32").Value = 1000

# Row 33
$ws.Range("D33").Value = 44348
$ws.Range("L33").Value = 'Segunda'
$ws.Range("M33").Value = 200
$ws.Range("N33").Value = 13000
$ws.Range("O33").Value = 13500
$ws.Range("P33").Value = 13250
$ws.Range("S33").Value = 883

# Row 34
$ws.Range("D34").Value = 44327
$ws.Range("L34").Value = 'Primera'
$ws.Range("M34").Value = 100
$ws.Range("N34").Value = 17000
$ws.Range("O34").Value = 17000
$ws.Range("P34").Value = 17000
$ws.Range("S34").Value = 1133

# Row 35
$ws.Range("D35").Value = 44327
$ws.Range("L35").Value = 'Segunda'
$ws.Range("M35").Value = 200
$ws.Range("N35").Value = 14000
$ws.Range("O35").Value = 14500
$ws.Range("P35").Value = 14250
$ws.Range("S35").Value = 950

# Row 36
$ws.Range("L36").Value = 'Primera'
$ws.Range("M36").Value = 100
$ws.Range("N36").Value = 17000
$ws.Range("O36").Value = 17000
$ws.Range("P36").Value = 17000
$ws.Range("S36").Value = 1133

# Row 37
$ws.Range("D37").Value = 44330
$ws.Range("L37").Value = 'Segunda'
$ws.Range("M37").Value = 200
$ws.Range("O37").Value = 14500
$ws.Range("P37").Value = 14250
$ws.Range("S37").Value = 950

# Row 38
$ws.Range("D38").Value = 44309

# Row 39
$ws.Range("D39").Value = 44309
$ws.Range("O39").Value = 14500
$ws.Range("S39").Value = 950
